$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit permutes the data of rows 67-73 (row 70 is left untouched) in two
# independent 3-cycles of "species observation" records:
#   67 <- 69 <- 71 <- 67   (new67 = old69, new69 = old71, new71 = old67)
#   68 <- 72 <- 73 <- 68   (new68 = old72, new72 = old73, new73 = old68)
# Columns A,B,D,E,F,G,H,P,Q,R carry the record's data and move with it.
# K ("Ålder-Stadium") and M ("Aktivitet") are mutually exclusive per record
# and also move with the record they belong to.

function Get-RowData($r) {
    $data = @{}
    $data.A = $ws.Range("A$r").Value2
    $data.B = $ws.Range("B$r").Value2
    $data.D = $ws.Range("D$r").Value2
    $data.E = $ws.Range("E$r").Value2
    $data.F = $ws.Range("F$r").Value2
    $data.G = $ws.Range("G$r").Value2
    $data.H = $ws.Range("H$r").Value2
    $data.K = $ws.Range("K$r").Value2
    $data.M = $ws.Range("M$r").Value2
    $data.P = $ws.Range("P$r").Value2
    $data.Q = $ws.Range("Q$r").Value2
    $data.R = $ws.Range("R$r").Value2
    return $data
}

# Snapshot the "before" state of every row involved, before any writes.
$row67 = Get-RowData 67
$row68 = Get-RowData 68
$row69 = Get-RowData 69
$row71 = Get-RowData 71
$row72 = Get-RowData 72
$row73 = Get-RowData 73

function Set-RowData($r, $data, $src) {
    $ws.Range("A$r").Value2 = $data.A
    $ws.Range("B$r").Value2 = $data.B
    $ws.Range("D$r").Value2 = $data.D
    $ws.Range("E$r").Value2 = $data.E
    $ws.Range("F$r").Value2 = $data.F
    $ws.Range("G$r").Value2 = $data.G
    $ws.Range("H$r").Value2 = $data.H

    # K / M: only touch them if this record's Age-Stage / Activity info is
    # actually changing, so untouched blank cells aren't disturbed.
    if ($data.K -ne $src.K) {
        if ($data.K -eq $null) {
            $ws.Range("K$r").ClearContents()
        } else {
            $ws.Range("K$r").Value2 = $data.K
        }
    }
    if ($data.M -ne $src.M) {
        if ($data.M -eq $null) {
            $ws.Range("M$r").ClearContents()
        } else {
            $ws.Range("M$r").Value2 = $data.M
        }
    }

    $ws.Range("P$r").Value2 = $data.P
    $ws.Range("Q$r").Value2 = $data.Q
    $ws.Range("R$r").Value2 = $data.R
}

# Apply the new content according to the two 3-cycles.
Set-RowData 67 $row69 $row67
Set-RowData 69 $row71 $row69
Set-RowData 71 $row67 $row71

Set-RowData 68 $row72 $row68
Set-RowData 72 $row73 $row72
Set-RowData 73 $row68 $row73
